# Add a new "2022-Q3" quarter:
#  - 总计 (summary) sheet: insert a new top data row for 2022-Q3 and
#    renumber the existing rows' index column.
#  - Insert a brand-new worksheet "2022-Q3" (holdings detail) right after
#    "总计", built from a copy of the existing "2022-Q2" sheet so it
#    inherits identical formatting, then overwritten with the new data.

$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $val) {
    # Force the cell to remain text (Excel would otherwise silently
    # coerce numeric-looking strings like "014185" or "3.62" into
    # numbers and lose the leading zero / fixed formatting).
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. "总计" sheet: insert new row 2 for 2022-Q3, shift the rest down.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()

# Re-apply the index-column (A) style from the row below, which still
# carries the original "s=2" bold/border formatting.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("B2").Style = "Normal"
$summary.Range("C2").Value = 18
$summary.Range("C2").Style = "Normal"
$summary.Range("D2").Value = 1.05
$summary.Range("D2").Style = "Normal"

# Renumber the index column (A) for the rows that shifted down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4

# ---------------------------------------------------------------------
# 2. New "2022-Q3" worksheet, positioned right after "总计".
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($null, $summary)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# The template ("2022-Q2") has 23 data rows (rows 2-24); the new sheet
# only needs 18 (rows 2-19) - clear the extra tail rows entirely.
$q3.Range("A20:H24").Clear()

$data = @(
    @(0, '014185', '招商专精特新股票A', '3.62', '81.87', '8.02', '0.2903', 2),
    @(1, '014186', '招商专精特新股票C', '2.63', '81.87', '8.02', '0.2109', 2),
    @(2, '162203', '泰达宏利稳定混合', '3.13', '91.16', '4.65', '0.1455', 4),
    @(3, '217013', '招商中小盘精选混合', '2.52', '80.73', '4.84', '0.1220', 6),
    @(4, '501030', '汇添富中证环境治理指数（LOF）A', '3.12', '92.74', '2.60', '0.0811', 3),
    @(5, '014320', '德邦半导体产业混合C', '1.38', '91.65', '4.14', '0.0571', 8),
    @(6, '164908', '交银施罗德中证环境治理指数（LOF）', '1.57', '93.62', '2.58', '0.0405', 4),
    @(7, '501031', '汇添富中证环境治理指数（LOF）C', '1.30', '92.74', '2.60', '0.0338', 3),
    @(8, '001531', '招商安益灵活配置混合', '0.57', '70.35', '3.71', '0.0211', 8),
    @(9, '014319', '德邦半导体产业混合A', '0.37', '91.65', '4.14', '0.0153', 8),
    @(10, '009719', '招商增浩一年定期开放混合C', '1.34', '23.60', '0.77', '0.0103', 8),
    @(11, '004352', '北信瑞丰研究精选股票', '0.45', '92.65', '1.29', '0.0058', 4),
    @(12, '002068', '东方多策略灵活配置混合C', '0.26', '55.14', '2.06', '0.0054', 10),
    @(13, '009718', '招商增浩一年定期开放混合A', '0.70', '23.60', '0.77', '0.0054', 8),
    @(14, '015641', '银华数字经济股票A', '0.14', '90.09', '2.70', '0.0038', 10),
    @(15, '013413', '交银施罗德中证环境治理指数（LOF）C', '0.09', '93.62', '2.58', '0.0023', 4),
    @(16, '400023', '东方多策略灵活配置混合A', '0.03', '55.14', '2.06', '0.0006', 10),
    @(17, '015642', '银华数字经济股票C', '0.01', '90.09', '2.70', '0.0003', 10)
)

$r = 2
foreach ($row in $data) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    Set-TextValue $q3.Cells.Item($r, 2) $row[1]
    Set-TextValue $q3.Cells.Item($r, 3) $row[2]
    Set-TextValue $q3.Cells.Item($r, 4) $row[3]
    Set-TextValue $q3.Cells.Item($r, 5) $row[4]
    Set-TextValue $q3.Cells.Item($r, 6) $row[5]
    Set-TextValue $q3.Cells.Item($r, 7) $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Restore "总计" as the active sheet/selection, matching the original file.
$summary.Activate()
